$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated column D values in rows 13-16 (the data table used to
# have the output column duplicated in C and D; now only C is kept).
$ws.Range("D13:D16").ClearContents()

# Update the view state to match the saved selection/scroll position.
$ws.Range("H19").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
